$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.212.37"
$ws.Range("E2").Value = "  +2.80%  "
$ws.Range("D3").Value = "3.631.71"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'197.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.98%  "
$ws.Range("D6").Value = "'576.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").Value = "3.627.49"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("D8").Value = "'0.619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.677"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +8.06%  "
$ws.Range("D12").Value = "'56.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.15%  "
$ws.Range("E13").Value = "  +17.31%  "
$ws.Range("D14").Value = "'10.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "4.218.93"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "3.639.14"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "'12.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("D19").Value = "68.199.49"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").Value = "'18.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("D22").Value = "'402.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.30%  "
$ws.Range("D23").Value = "'13.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +30.33%  "
$ws.Range("D24").Value = "'4.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'85.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").Value = "'2.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.24%  "
$ws.Range("D27").Value = "'12.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("D28").Value = "'3.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.75%  "
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").Value = "'8.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +23.49%  "
$ws.Range("D31").Value = "'9.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.77%  "
$ws.Range("D32").Value = "'31.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'690.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +17.14%  "
$ws.Range("D34").Value = "'12.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("E35").Value = "  +5.32%  "
$ws.Range("D36").Value = "'64.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").Value = "'42.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.29%  "
$ws.Range("D38").Value = "'0.428"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.06%  "
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("D40").Value = "0.0₃0788"
$ws.Range("E40").Value = "  +7.95%  "
$ws.Range("E41").Value = "  +8.28%  "
$ws.Range("E42").Value = "  +21.73%  "
$ws.Range("D43").Value = "'3.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.56%  "
$ws.Range("D44").Value = "3.222.42"
$ws.Range("E44").Value = "  +16.41%  "
$ws.Range("D45").Value = "'3.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +43.66%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'0.0421"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.80%  "
$ws.Range("D48").Value = "'8.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.79%  "
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").Value = "'3.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'2.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.33%  "
